# Update "days since" (column F) counters across sheets to reflect the
# regenerated output date (gh-pages rebuild at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 14595
$ws1.Range("F4").Value = 14823
$ws1.Range("F13").Value = 1595
$ws1.Range("F18").Value = 934
$ws1.Range("F23").Value = 3491
$ws1.Range("F27").Value = 643
$ws1.Range("F30").Value = 1872
$ws1.Range("F32").Value = 1527
$ws1.Range("F33").Value = 328
$ws1.Range("F35").Value = 7226
$ws1.Range("F36").Value = 5062
$ws1.Range("F39").Value = 711
$ws1.Range("F40").Value = 3353
$ws1.Range("F43").Value = 355
$ws1.Range("F44").Value = 141
$ws1.Range("F46").Value = 4478
$ws1.Range("F47").Value = 666

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 7855
$ws3.Range("F3").Value = 293
$ws3.Range("F4").Value = 1017

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7855
$ws4.Range("F4").Value = 293
$ws4.Range("F5").Value = 1017
$ws4.Range("F7").Value = 14595
$ws4.Range("F8").Value = 14823
$ws4.Range("F17").Value = 1595
$ws4.Range("F22").Value = 3491
$ws4.Range("F25").Value = 643
$ws4.Range("F27").Value = 1872
$ws4.Range("F35").Value = 7226
$ws4.Range("F36").Value = 5062
$ws4.Range("F39").Value = 3353
$ws4.Range("F41").Value = 355
$ws4.Range("F42").Value = 141
$ws4.Range("F45").Value = 666
